$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add IS Devices")

# Replace the label strings in E8:E15 ("1 / 1" ... "1 / 8") with the plain
# numeric value 1, per the updated test data implementation.
for ($r = 8; $r -le 15; $r++) {
    $ws.Cells.Item($r, 5).Value = 1
}

# Setting a numeric value into these cells (previously text with a quote
# prefix / style index 11) re-derives their cell style. Re-apply the
# original formatting (style index 11) by pasting formats from a donor
# cell elsewhere in the workbook that still carries that exact style, so
# the cells keep their original look while now holding numeric values.
$donor = $wb.Worksheets.Item("Sheet1").Cells.Item(4, 8)
$donor.Copy()
$ws.Range("E8:E15").PasteSpecial(-4122)  # xlPasteFormats

# Update the active selection on this sheet to E9:E15 (activeCell E9)
$ws.Range("E9:E15").Select()
